$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for rows 99-110 (Trial 96-107), columns A (1) through AI (35)
$rows = @(
  @('Trial    96', 1465.93, 10.8856, 17, 138.1, 1.1, 20, 45.6666, 2.2, 11, 19.9, 1.1, 8, 5.76666, 2.1, 5, 20.1666, 48.3666, 3, 2.1, 3.46666, 3, 17.1, 71.7333, 2, 20.3, 35.1333, 3, 2.66666, 103.033, 4, 2.93333, 9.76666, 12.1581, 49483.4),
  @('Trial    97', 1654.02, 14.5047, 26, 119.6, 0, 19, 28.5333, 3.73333, 13, 11.5333, 0, 12, 14.8, 2.36666, 11, 15.3, 2.83333, 5, 14.0667, 3.9, 4, 17.2333, 46.5666, 3, 13.2333, 10.2, 3, 1.9, 73.2999, 5, 1.06667, 14.8667, 13.2942, 46117.5),
  @('Trial    98', 1557.9, 11.6232, 24, 138.833, 3.4, 22, 46.0666, 3.8, 12, 11, 29.2, 11, 9.36666, 3.4, 4, 25.6666, 4.03333, 4, 11.2, 5.29999, 5, 16.4333, 8.59999, 1, 7.86666, 136.467, 4, 4.39999, 28.4333, 8, 5.13333, 29.0666, 11.6261, 47213.5),
  @('Trial    99', 2476.81, 13.9381, 28, 182.5, 1.53333, 28, 33.2333, 5.13333, 14, 20.2666, 1.53333, 10, 25.0666, 2.1, 6, 22.1666, 2.3, 8, 14.1, 6.96666, 13, 44.3666, 6.19999, 2, 13.2667, 28.3666, 3, 2.13333, 110.7, 5, 5.89999, 35.7666, 15.8732, 85223.4),
  @('Trial   100', 1045.09, 11.3844, 19, 94.6332, 2.3, 13, 25.6, 3.56666, 9, 14.2, 2.3, 8, 14.4333, 3.03333, 1, 0.133333, 3.4, 4, 11.0667, 47.1333, 5, 17.2333, 46.1666, 1, 8.16666, 81.3333, 1, 0.666666, 25.1666, 3, 2, 10.3667, 14.6768, 40757.5),
  @('Trial   101', 476.747, 13.3543, 8, 37.1666, 0.933332, 6, 4.23333, 7.29999, 4, 6.29999, 0.933332, 7, 9.76666, 6.36666, 2, 1.76666, 6.83333, 1, 7.59999, 27.2666, 1, 1.3, 7.66666, 1, 5.76666, 35.6666, 0, 0, '-', 0, 0, '-', 19.0708, 20653.7),
  @('Trial   102', 821.037, 13.7297, 14, 62.0333, 4.43333, 9, 15.6667, 4.69999, 9, 7.43332, 12.0333, 7, 8.39999, 12.4, 5, 15.0333, 4.43333, 1, 7.89999, 59.5999, 1, 0.966665, 67.8999, 1, 3.56666, 6.53333, 0, 0, '-', 3, 2.7, 10.7667, 14.4364, 26158.8),
  @('Trial   103', 1090.99, 10.2249, 28, 110.9, 0.366666, 16, 34.4333, 2.66666, 10, 15.4, 13.1, 12, 17.5, 0.366666, 9, 18.4666, 1.86666, 4, 6.23333, 2.93333, 2, 1.8, 9.26666, 1, 14.5333, 99.1999, 0, 0, '-', 4, 1.46667, 12.5, 14.2485, 46094.1),
  @('Trial   104', 417.726, 9.66214, 1, 43.3666, 0, 5, 5.2, 4.43333, 1, 4.43333, 0, 1, 6.26666, 36.2333, 1, 6.59999, 27.3666, 1, 11.3667, 15.5, 0, 0, '-', 1, 9.29999, 5.39999, 0, 0, '-', 0, 0, '-', 20.6066, 26747.4),
  @('Trial   105', 940.208, 14.4722, 1, 65.0999, 0, 14, 15.1667, 0, 1, 6.19999, 43.8, 3, 10.1, 2.4, 1, 1.76666, 7.29999, 2, 9.43332, 17.4333, 2, 8.66666, 0.033333, 2, 7.69999, 56.2666, 1, 3.6, 33.0333, 1, 2, 10.3, 16.2782, 31742.5),
  @('Trial   106', 348.143, 9.05052, 1, 38.5, 0, 5, 5.23333, 2.16666, 1, 6.59999, 2.76666, 1, 7.29999, 10.0333, 0, 0, '-', 1, 7.36666, 19.1, 0, 0, '-', 1, 9.63332, 28.3, 0, 0, '-', 1, 2.13333, 0, 18.7527, 21659.4),
  @('Trial   107', 628.822, 13.2663, 1, 47.4333, 0.099999, 10, 3.96666, 0.099999, 1, 7.36666, 6.29999, 1, 6.43333, 14, 1, 2.8, 20.7666, 1, 7.49999, 23.9333, 2, 4.43333, 0.133333, 1, 11.6, 35.3666, 1, 1.46667, 2.2, 1, 1.73333, 4.16666, 21.5368, 30646.8)
)

$startRow = 99
for ($i = 0; $i -lt $rows.Length; $i++) {
  $rowData = $rows[$i]
  $r = $startRow + $i
  for ($j = 0; $j -lt $rowData.Length; $j++) {
    $ws.Cells.Item($r, $j + 1).Value = $rowData[$j]
  }
}
